# CobaltUsers.xlsx update
# - Rename Sheet2 -> Emails, populate it with an Email/Password table
# - Append 28 new test users (rows 53-80) to the Users sheet, each with
#   a "do not use" note, Locked=N, and a mailto hyperlink in column G
# - Add a handful of blank bordered rows below (81-84) matching the
#   original author's leftover formatting

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet2 -> Emails
# ---------------------------------------------------------------------
$emailsWs = $wb.Worksheets.Item("Sheet2")
$emailsWs.Name = "Emails"

$emailsWs.Range("A1").Value = "Email"
$emailsWs.Range("B1").Value = "Password"
$emailsWs.Range("A2").Value = "tr-anz-tester1@yandex.com"
$emailsWs.Range("B2").Value = "tranztest"
$emailsWs.Range("A3").Value = "tr-anz-tester2@yandex.com"
$emailsWs.Range("B3").Value = "tranztest"

$emailsWs.Columns.Item(1).ColumnWidth = 26.28515625
$emailsWs.Columns.Item(2).ColumnWidth = 13.85546875

$emailsWs.Range("A1:B3").Select()

# ---------------------------------------------------------------------
# Users sheet -> append new test accounts
# ---------------------------------------------------------------------
$usersWs = $wb.Worksheets.Item("Users")

$newUsers = @(
    @{ Name = "SearchOpenWebUser1"; Email = "SearchOpenWeb@mailinator.com " },
    @{ Name = "FFHUser1";           Email = "FFHUser1@mailinator.com " },
    @{ Name = "FFHUser2";           Email = "FFHUser2@mailinator.com" },
    @{ Name = "FFHUser3";           Email = "FFHUser3@mailinator.com" },
    @{ Name = "FFHUser4";           Email = "FFHUser4@mailinator.com" },
    @{ Name = "FrontEndUser1";      Email = "FrontEndUser1@mailinator.com" },
    @{ Name = "FrontEndUser2";      Email = "FrontEndUser2@mailinator.com" },
    @{ Name = "FrontEndUser3";      Email = "FrontEndUser3@mailinator.com" },
    @{ Name = "FrontEndUser4";      Email = "FrontEndUser4@mailinator.com" },
    @{ Name = "FrontEndUser5";      Email = "FrontEndUser5@mailinator.com" },
    @{ Name = "FrontEndUser6";      Email = "FrontEndUser6@mailinator.com" },
    @{ Name = "FrontEndUser7";      Email = "FrontEndUser7@mailinator.com" },
    @{ Name = "FrontEndUser8";      Email = "FrontEndUser8@mailinator.com" },
    @{ Name = "FrontEndUser9";      Email = "FrontEndUser9@mailinator.com" },
    @{ Name = "FrontEndUser10";     Email = "FrontEndUser10@mailinator.com" },
    @{ Name = "UrlUser1";           Email = "UrlUser1@mailinator.com" },
    @{ Name = "UrlUser2";           Email = "UrlUser2@mailinator.com" },
    @{ Name = "UrlUser3";           Email = "UrlUser3@mailinator.com" },
    @{ Name = "LinkingUser1";       Email = "LinkingUser1@mailinator.com " },
    @{ Name = "LoginUser1";         Email = "LoginUser1@mailinator.com " },
    @{ Name = "LoginUser2";         Email = "LoginUser2@mailinator.com " },
    @{ Name = "LoginUser3";         Email = "LoginUser3@mailinator.com " },
    @{ Name = "LoginUser4";         Email = "LoginUser4@mailinator.com " },
    @{ Name = "LoginUser5";         Email = "LoginUser5@mailinator.com " },
    @{ Name = "LoginUser6";         Email = "LoginUser6@mailinator.com " },
    @{ Name = "LoginUser7";         Email = "LoginUser7@mailinator.com " },
    @{ Name = "CpetUser1";          Email = "CpetUser1@mailinator.com " },
    @{ Name = "CpetUser2";          Email = "CpetUser2@mailinator.com " }
)

$row = 53
foreach ($u in $newUsers) {
    $usersWs.Range("A$row").Value = $u.Name
    $usersWs.Range("B$row").Value = "Password1"

    $eCell = $usersWs.Range("E$row")
    $eCell.Value = "THIS IS IN USE 24/7 - DO NOT USE!"
    $eCell.Borders.Item(7).LineStyle = 1
    $eCell.Borders.Item(10).LineStyle = 1

    $fCell = $usersWs.Range("F$row")
    $fCell.Value = "N"
    $fCell.Borders.Item(7).LineStyle = 1
    $fCell.Borders.Item(10).LineStyle = 1

    $gCell = $usersWs.Range("G$row")
    $gCell.Value = $u.Email
    $usersWs.Hyperlinks.Add($gCell, "mailto:" + $u.Email.Trim())

    $row = $row + 1
}

# Trailing formatted-but-empty rows (81-84), matching the source file
for ($r = 81; $r -le 84; $r++) {
    $eCell = $usersWs.Range("E$r")
    $eCell.Borders.Item(7).LineStyle = 1
    $eCell.Borders.Item(10).LineStyle = 1
}

$usersWs.Columns.Item(1).ColumnWidth = 21.140625
$usersWs.Columns.Item(5).ColumnWidth = 39.28515625

$usersWs.Range("C80").Select()
